# "Fruta / hortaliza, semanal" — weekly price-sheet update.
#
# A new weekly observation is inserted as a new data row right before the
# current row 206 (pushing the existing rows 206-288 down to 207-289), for
# the "Feria Lagunitas de Puerto Montt - Cilantro" series. All the other
# (already-shifted) rows keep their original data untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 206; everything from the old row 206 down shifts to
# row 207 onward, and the sheet's used range grows from R288 to R289.
$ws.Rows.Item(206).EntireRow.Insert()

# Populate the newly inserted row 206 with the new weekly record.
$newRow = New-Object 'object[,]' 1,18
$newRow[0,0]  = 4                                        # Mercado ID
$newRow[0,1]  = "Feria Lagunitas de Puerto Montt"         # Mercado
$newRow[0,2]  = "Los Lagos"                               # Region
$newRow[0,3]  = 44704                                     # Fecha
$newRow[0,4]  = 10                                        # Codreg
$newRow[0,5]  = 100112040                                 # Categoria ID
$newRow[0,6]  = "Cilantro"                                # Categoria
$newRow[0,7]  = "Sin especificar"                         # Variedad
$newRow[0,8]  = "Primera"                                 # Calidad
$newRow[0,9]  = 70                                        # Volumen
$newRow[0,10] = 6000                                      # Precio minimo
$newRow[0,11] = 6000                                      # Precio maximo
$newRow[0,12] = 6000                                      # Precio promedio ponderado
$newRow[0,13] = "`$/docena de atados (2 kilos)"            # Unidad de comercializacion
$newRow[0,14] = "Región de La Araucanía"                  # Origen
$newRow[0,15] = 3000                                      # Precio $/Kg
$newRow[0,16] = 2                                         # Kg o Unidades
$newRow[0,17] = "Hortaliza"                                # Clasificacion

$ws.Range("A206:R206").Value = $newRow

# Keep the date cell formatted like the rest of column D (Insert already
# copies formatting from the row above, but set it explicitly to be safe).
$ws.Cells.Item(206, 4).NumberFormat = $ws.Cells.Item(207, 4).NumberFormat
